$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 41.43823833333334
$ws.Range("H2").Value = 124.314715
$ws.Range("I2").Value = 0.981992391336623
$ws.Range("J2").Value = 0.9819923913366232
$ws.Range("M2").Value = 48.42420966666666
$ws.Range("N2").Value = 145.272629
$ws.Range("O2").Value = 0.6311762527593259
$ws.Range("P2").Value = 0.6311762527593258
$ws.Range("Q2").Value = 2006.613941270637
$ws.Range("R2").Value = 18059.52547143574
$ws.Range("S2").Value = 0.6198102778020191
$ws.Range("T2").Value = 0.6198102778020192

$ws.Range("G3").Value = 41.43823833333334
$ws.Range("H3").Value = 124.314715
$ws.Range("I3").Value = 0.981992391336623
$ws.Range("J3").Value = 0.9819923913366232
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("O3").Value = 0.08928392431779728
$ws.Range("P3").Value = 0.08928392431779726
$ws.Range("Q3").Value = 283.8483965203289
$ws.Range("R3").Value = 2554.63556868296
$ws.Range("S3").Value = 0.08767613434875181
$ws.Range("T3").Value = 0.08767613434875182

$ws.Range("G4").Value = 41.43823833333334
$ws.Range("H4").Value = 124.314715
$ws.Range("I4").Value = 0.981992391336623
$ws.Range("J4").Value = 0.9819923913366232
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2795398229228769
$ws.Range("P4").Value = 0.2795398229228769
$ws.Range("Q4").Value = 888.7034380098232
$ws.Range("R4").Value = 7998.33094208841
$ws.Range("S4").Value = 0.274505979185852
$ws.Range("T4").Value = 0.2745059791858521

$ws.Range("I5").Value = 0.006845967574057415
$ws.Range("J5").Value = 0.006845967574057417
$ws.Range("M5").Value = 48.42420966666666
$ws.Range("N5").Value = 145.272629
$ws.Range("O5").Value = 0.6311762527593259
$ws.Range("P5").Value = 0.6311762527593258
$ws.Range("Q5").Value = 13.98912465797433
$ws.Range("R5").Value = 125.902121921769
$ws.Range("S5").Value = 0.004321012159905412
$ws.Range("T5").Value = 0.004321012159905413

$ws.Range("I6").Value = 0.006845967574057415
$ws.Range("J6").Value = 0.006845967574057417
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("O6").Value = 0.08928392431779728
$ws.Range("P6").Value = 0.08928392431779726
$ws.Range("S6").Value = 0.0006112348507642364
$ws.Range("T6").Value = 0.0006112348507642365

$ws.Range("I7").Value = 0.006845967574057415
$ws.Range("J7").Value = 0.006845967574057417
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2795398229228769
$ws.Range("P7").Value = 0.2795398229228769
$ws.Range("Q7").Value = 6.195602912245999
$ws.Range("R7").Value = 55.76042621021399
$ws.Range("S7").Value = 0.001913720563387767
$ws.Range("T7").Value = 0.001913720563387767

$ws.Range("I8").Value = 0.01116164108931947
$ws.Range("J8").Value = 0.01116164108931947
$ws.Range("M8").Value = 48.42420966666666
$ws.Range("N8").Value = 145.272629
$ws.Range("O8").Value = 0.6311762527593259
$ws.Range("P8").Value = 0.6311762527593258
$ws.Range("Q8").Value = 22.80781889440322
$ws.Range("R8").Value = 205.270370049629
$ws.Range("S8").Value = 0.007044962797401183
$ws.Range("T8").Value = 0.007044962797401184

$ws.Range("I9").Value = 0.01116164108931947
$ws.Range("J9").Value = 0.01116164108931947
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("O9").Value = 0.08928392431779728
$ws.Range("P9").Value = 0.08928392431779726
$ws.Range("Q9").Value = 3.226312091304889
$ws.Range("S9").Value = 0.0009965551182812156
$ws.Range("T9").Value = 0.000996555118281216

$ws.Range("I10").Value = 0.01116164108931947
$ws.Range("J10").Value = 0.01116164108931947
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2795398229228769
$ws.Range("P10").Value = 0.2795398229228769
$ws.Range("R10").Value = 90.91159980137398
$ws.Range("S10").Value = 0.003120123173637071
$ws.Range("T10").Value = 0.003120123173637072
